# issue #5: stock data output to json file
#
# On the 股票 (stock) sheet, insert a new "property_category" column
# (constant value "stock") between "total" and "date", pushing
# date / legislator_name / legislator_id one column to the right
# (H->I, I->J, J->K). Also fixes a stray " '" typo on the "仁寶" stock name.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(6)

$xlPasteFormats = -4122
$xlPasteValues = -4163

$lastRow = 11
$headerStyleCell = $ws.Range("B1")   # bold header style, used as a format template

function Set-TextValue($cell, [string]$text) {
    # Writing a date-shaped literal straight into .Value lets the host
    # auto-coerce it into a real Date; routing it through a formula and
    # converting the formula result to a literal value keeps it a plain
    # string. Done in place so any pre-existing cell style survives.
    $cell.Formula = '="' + $text + '"'
    $cell.Copy()
    $cell.PasteSpecial($xlPasteValues)
}

for ($r = 1; $r -le $lastRow; $r++) {
    $hCell = $ws.Cells.Item($r, 8)   # old date / "date" header
    $iCell = $ws.Cells.Item($r, 9)   # old legislator_name / header
    $jCell = $ws.Cells.Item($r, 10)  # old legislator_id / header
    $kCell = $ws.Cells.Item($r, 11)  # brand-new column

    $hVal = $hCell.Value()
    $iVal = $iCell.Value()
    $jVal = $jCell.Value()

    if ($r -eq 1) {
        # ---- header row: everything is text ----
        $kCell.Value = $jVal
        $headerStyleCell.Copy()
        $kCell.PasteSpecial($xlPasteFormats)

        Set-TextValue $jCell $iVal
        Set-TextValue $iCell $hVal
        Set-TextValue $hCell "property_category"
    } else {
        # ---- data rows ----
        $kCell.Value = $jVal          # legislator_id (numeric) -> K

        Set-TextValue $jCell $iVal    # legislator_name -> J
        Set-TextValue $iCell $hVal    # date -> I
        Set-TextValue $hCell "stock"  # new property_category value -> H
    }
}

# ---- Typo fix: "仁寶 '" -> "仁寶" (row 10, column B) ----
$ws.Cells.Item(10, 2).Value = "仁寶"
